# Finished authorization and added users link in navbar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Developmnet PB")   # "Development PB" sheet (the active/first tab)

# Row 49 gains a Status (column C) value matching the "Not Started" entries above.
$ws.Range("C49").Value = "Not Started"

# New row 50: "Add ""active"" attribute in the bootstrap nav panel"
$ws.Range("A50").Value = 'Add "active" attribute in the bootstrap nav panel'
$ws.Range("B50").Interior.Color = $ws.Range("B45").Interior.Color
$ws.Range("C50").Value = "Not Started"

# New row 51: "Make ""Patient Handling System"" look better in the nav bar"
$ws.Range("A51").Value = 'Make "Patient Handling System" look better in the nav bar'
$ws.Range("B51").Interior.Color = $ws.Range("B45").Interior.Color
$ws.Range("C51").Value = "Not Started"

$null = $ws.Range("E49").Select()
